# refactor: rename column Transportadora to Importadora
#
# The "Transportadora" header lives in cell C1 of the (only) worksheet.
# Update its text in place, the same way a user would by typing a new
# value into the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 3).Value = "Importadora"
